$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / top rows -------------------------------------------------
# B1 "Quantidade" label (shared-string content unchanged, just re-asserted)
$ws.Range("B1").Value = "Quantidade"

# A2:A11 become the new (much shorter) product list. These rows also pick
# up the "0.00" number format (matches the target style used for the new
# list) in addition to keeping their text content.
$ws.Range("A2").Value = "Bromoprida"
$ws.Range("A3").Value = "Deflazacort"
$ws.Range("A4").Value = "Dustasterida"
$ws.Range("A5").Value = "Terbinafina"
$ws.Range("A6").Value = "SAME"
$ws.Range("A7").Value = "Passiflora"
$ws.Range("A8").Value = "Açafrão"
$ws.Range("A9").Value = "Vit D"
$ws.Range("A10").Value = "Carbotil UG"
$ws.Range("A11").Value = "Cápsula 00 Incolor"
$ws.Range("A2:A11").NumberFormat = "0.00"

# --- Remove the old (now obsolete) product rows -------------------------
# Rows 12-114 used to carry the old, much longer product catalogue; the
# new sheet just leaves column A blank for them.
For ($r = 12; $r -le 114; $r++) {
    $ws.Cells.Item($r, 1).Value = ""
}

# --- New helper cell with a plain (non bold/italic/underline) Arial font
$d16 = $ws.Range("D16")
$d16.Font.Name = "Arial"
$d16.Font.Size = 11
$d16.Font.Bold = $false
$d16.Font.Italic = $false
$d16.Font.Underline = $false

# --- Selection / view state ---------------------------------------------
$ws.Range("A2").Select()
